# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (column G) of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
